# Update Name of Algo
# Updates imputed values in columns C and D (KNN result data) for a set of rows
# to reflect the re-run/updated algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C11").Value = -11.515
$ws.Range("D11").Value = -7.436999999999999
$ws.Range("C12").Value = -10.534
$ws.Range("C15").Value = -13.474
$ws.Range("D23").Value = -8.247000000000002
$ws.Range("C27").Value = -13.148
$ws.Range("C28").Value = -12.302
$ws.Range("D28").Value = -7.801
$ws.Range("C31").Value = -12.987
$ws.Range("C32").Value = -12.087
$ws.Range("D32").Value = -7.540999999999999
$ws.Range("D34").Value = -7.933
$ws.Range("C36").Value = -12.622
$ws.Range("D36").Value = -8.122
$ws.Range("D37").Value = -8.097999999999999
$ws.Range("C38").Value = -12.577
$ws.Range("D42").Value = -8.286999999999999
$ws.Range("C46").Value = -13.836
$ws.Range("D49").Value = -8.289999999999999
$ws.Range("C54").Value = -12.528
$ws.Range("D54").Value = -7.531000000000001
$ws.Range("C55").Value = -13.868
$ws.Range("C56").Value = -13.593
$ws.Range("C67").Value = -11.682
$ws.Range("C69").Value = -10.958
$ws.Range("C72").Value = -11.753
$ws.Range("C73").Value = -12.446
$ws.Range("D78").Value = -8.236000000000001
$ws.Range("D80").Value = -8.050000000000001
$ws.Range("C83").Value = -13.465
$ws.Range("C86").Value = -14.085
$ws.Range("C91").Value = -12.173
$ws.Range("C93").Value = -10.844
$ws.Range("D97").Value = -7.814
$ws.Range("C99").Value = -11.242
$ws.Range("D99").Value = -7.428
$ws.Range("D100").Value = -7.706999999999999
$ws.Range("D101").Value = -7.827000000000001
$ws.Range("C104").Value = -12.573
$ws.Range("C105").Value = -13.389
